$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.01725120502155203, 1.384110966838059, 8.376137918165924, 2.894155821334768, 2.959148581664683, 23)
    3  = @(-0.1059444359486, 1.363387133547353, 7.96756717039053, 2.822687933582197, 2.887077487796805, 22)
    4  = @(-0.5314575476434735, 0.935032954889049, 3.962372580821567, 1.990570918310013, 1.965686070731558, 21)
    5  = @(-0.1228908575652036, 0.5948419731728133, 0.8938785104026519, 0.9454514849544908, 0.9617836213211481, 20)
    6  = @(-0.07803016097268373, 0.595167622792783, 0.6624405579447651, 0.8139045140216173, 0.8323555993064871, 19)
    7  = @(-0.1177682127167521, 0.5278966303169981, 0.6234154974736094, 0.7895666516980118, 0.8033690317453926, 18)
    8  = @(0.008373115946536357, 0.403612244087397, 0.3646454851063948, 0.6038588287889768, 0.6223835931625885, 17)
    9  = @(0.005946962572950867, 0.3740864020203869, 0.2564549243757536, 0.5064137877030538, 0.5229858460921323, 16)
    10 = @(0.01722052028762557, 0.3408956206005759, 0.2330792879415247, 0.4827828579615527, 0.4994097318934355, 15)
    11 = @(0.04981916865326922, 0.3023652478618333, 0.214504976889282, 0.4631468200142175, 0.4778414914637334, 14)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
